# Refresh the cryptos price/volume snapshot (GitHub Actions style update).
#
# The "Price" column stores values as TEXT (e.g. "25.903.91", "1.007") rather
# than numbers, since many of these strings are not valid numeric literals
# (multiple "." separators) and the ones that are numeric-looking must keep
# trailing zeros ("1.010" must not collapse to "1.01"). Assigning a plain
# numeric-looking string to Range.Value makes Excel auto-convert it to a
# number, so for those cells we use the standard Excel "treat as text"
# idiom: a leading apostrophe (quote-prefix), exactly as typing '1.007 into a
# cell in the Excel UI would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.909.87'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.648.58'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').Value = "'215.67"
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = "'0.5090"
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').Value = "'1.007"
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = "'0.2576"
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.06428"
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').Value = "'19.74"
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').Value = "'0.07784"
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.320"
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.634.09'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = "'0.5480"
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '0.0₅7913'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = "'65.26"
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = '26.014.71'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = "'1.008"
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').Value = "'197.70"
$ws.Range('D20').Value = "'4.439"
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').Value = "'10.05"
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = "'6.069"
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('D23').Value = "'1.010"
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').Value = "'1.858"
$ws.Range('E24').Value = '  -2.97%  '
$ws.Range('D25').Value = "'141.52"
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = "'0.1148"
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = "'6.916"
$ws.Range('E27').Value = '  +3.04%  '
$ws.Range('D28').Value = "'15.76"
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = "'1.245"
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = "'0.05038"
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = "'3.281"
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').Value = "'3.211"
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('D33').Value = "'1.547"
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').Value = "'2.373"
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D35').Value = "'0.8950"
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = "'2.599"
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').Value = "'0.5558"
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').Value = '1.135.15'
$ws.Range('E38').Value = '  -3.48%  '
$ws.Range('D39').Value = "'0.01567"
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = "'1.010"
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('D41').Value = "'5.677"
$ws.Range('D42').Value = "'0.8171"
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').Value = "'99.93"
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('E44').Value = '  +6.41%  '
$ws.Range('D45').Value = '1.785.53'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D46').Value = "'0.4539"
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'55.39"
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = "'1.007"
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('D50').Value = "'1.008"
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('D51').Value = "'0.09577"
$ws.Range('E51').Value = '  +3.12%  '
